# The presentation's two theme parts (ppt/theme/theme1.xml = "Office Theme",
# ppt/theme/theme2.xml = "Integral") had their contents swapped: the slide
# master (which is wired to theme2.xml) now carries the Office Theme color
# scheme, while the Integral color scheme moves to theme1.xml (used by the
# notes master). The reachable PowerPoint object model surface for this is
# the Design/Master color scheme (12 theme colors, ppColorSchemeIndex 1-12:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) -- so re-point the slide
# master's 12 scheme colors at the "Office" palette that used to live in
# theme1.xml.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
